# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled
# update) — for every coin row on the sheet, write the newly scraped
# Price (column D) and Volume(1h) (column E) text, and for the two row
# pairs whose rank flipped (BitcoinCash <-> InternetComputer(DFINITY) at
# rows 28/29, and Quant <-> EnergySwap at rows 46/47) also rewrite the
# Coin name and Link columns (B/C) to match the new ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    # Force the value to be stored as text, matching the source
    # workbook where these cells are inline strings (prices use
    # dotted/locale formatting that must not be reinterpreted as
    # a number). Revert the cell style afterwards so no stray
    # number-format style lingers on the cell.
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '27.761.59'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = '1.855.39'
$ws.Range('E3').Value = '  +0.39%  '
$ws.Range('E4').Value = '  +0.22%  '
Set-TextCell 'D5' '312.22'
$ws.Range('E5').Value = '  -0.82%  '
$ws.Range('E6').Value = '  +0.22%  '
Set-TextCell 'D7' '0.4264'
$ws.Range('E7').Value = '  +0.40%  '
Set-TextCell 'D8' '0.3595'
$ws.Range('E8').Value = '  -1.66%  '
Set-TextCell 'D9' '0.07299'
$ws.Range('E9').Value = '  +0.01%  '
Set-TextCell 'D10' '0.8796'
$ws.Range('E10').Value = '  -1.55%  '
Set-TextCell 'D11' '20.81'
$ws.Range('E11').Value = '  -0.38%  '
$ws.Range('D12').Value = '1.879.82'
$ws.Range('E12').Value = '  +5.09%  '
Set-TextCell 'D13' '6.558'
$ws.Range('E13').Value = '  -0.37%  '
Set-TextCell 'D14' '5.347'
$ws.Range('E14').Value = '  +0.10%  '
Set-TextCell 'D15' '0.07013'
$ws.Range('E15').Value = '  +1.17%  '
$ws.Range('E16').Value = '  +0.24%  '
Set-TextCell 'D17' '79.88'
$ws.Range('E17').Value = '  +0.93%  '
Set-TextCell 'D18' '0.000008966'
$ws.Range('E18').Value = '  +0.75%  '
$ws.Range('E19').Value = '  +0.24%  '
Set-TextCell 'D20' '15.29'
$ws.Range('E20').Value = '  -1.21%  '
$ws.Range('D21').Value = '27.807.05'
$ws.Range('E21').Value = '  +0.23%  '
Set-TextCell 'D22' '5.016'
$ws.Range('E22').Value = '  +0.62%  '
Set-TextCell 'D23' '10.43'
$ws.Range('E23').Value = '  -2.03%  '
$ws.Range('D24').Value = '2.099.25'
$ws.Range('E24').Value = '  +3.11%  '
Set-TextCell 'D25' '1.981'
$ws.Range('E25').Value = '  +3.58%  '
Set-TextCell 'D26' '154.29'
$ws.Range('E26').Value = '  +0.24%  '
Set-TextCell 'D27' '18.52'
$ws.Range('E27').Value = '  -2.07%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 'D28' '5.280'
$ws.Range('E28').Value = '  -0.16%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell 'D29' '119.68'
$ws.Range('E29').Value = '  -2.37%  '
Set-TextCell 'D30' '1.901'
$ws.Range('E30').Value = '  -0.26%  '
Set-TextCell 'D31' '0.08895'
$ws.Range('E31').Value = '  -0.44%  '
Set-TextCell 'D32' '0.7604'
$ws.Range('E32').Value = '  -2.03%  '
Set-TextCell 'D33' '2.972'
$ws.Range('E33').Value = '  +0.32%  '
Set-TextCell 'D34' '4.526'
$ws.Range('E34').Value = '  -1.34%  '
Set-TextCell 'D35' '1.127'
$ws.Range('E35').Value = '  +1.93%  '
$ws.Range('E36').Value = '  +0.32%  '
Set-TextCell 'D37' '1.107'
$ws.Range('E37').Value = '  -0.23%  '
Set-TextCell 'D38' '0.05423'
$ws.Range('E38').Value = '  +0.63%  '
Set-TextCell 'D39' '0.01933'
$ws.Range('E39').Value = '  -0.98%  '
Set-TextCell 'D40' '2.821'
$ws.Range('E40').Value = '  +0.13%  '
Set-TextCell 'D41' '0.1676'
$ws.Range('E41').Value = '  +0.75%  '
Set-TextCell 'D42' '0.5100'
$ws.Range('E42').Value = '  -0.10%  '
Set-TextCell 'D43' '6.613'
$ws.Range('E43').Value = '  -4.20%  '
Set-TextCell 'D44' '8.449'
$ws.Range('E44').Value = '  +1.53%  '
Set-TextCell 'D45' '0.06533'
$ws.Range('E45').Value = '  -1.02%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D46' '10.36'
$ws.Range('E46').Value = '  -1.22%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 'D47' '105.64'
$ws.Range('E47').Value = '  +0.62%  '
Set-TextCell 'D48' '0.4677'
Set-TextCell 'D49' '0.9997'
$ws.Range('E49').Value = '  +0.26%  '
Set-TextCell 'D50' '1.623'
$ws.Range('E50').Value = '  -0.90%  '
Set-TextCell 'D51' '1.846'
$ws.Range('E51').Value = '  +3.76%  '
